$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the publishedDate column keeps storing plain text (as in the
# original workbook) instead of being auto-converted to numbers/dates.
$ws.Range("D4:D6").NumberFormat = "@"

# Row 4: now holds what used to be row 5's book (La comunidad del anillo)
$ws.Range("A4").Value = "Av6RMQEACAAJ"
$ws.Range("B4").Value = "La comunidad del anillo"
$ws.Range("C4").Value = "En la adormecida e idilica Comarca, un joven hobbit recibe un encargo : custodiar el Anillo Único y emprender el viaje para su destrucción en las Grietas del Destino. Acompañado por magos, hombres, elfos y enanos, atravesará la Tierra Media y se internará en las sombras de Mordor, perseguido siempre por las huestes de Sauron, el Señor Oscuro, dispuesto a recuperar su creación para establecer el dominio definitivo del Mal. (Source : 4e de couverture)."
$ws.Range("D4").Value = "2012"
$ws.Range("E4").Value = "John Ronald Reuel Tolkien"

# Row 5: now holds what used to be row 6's book (edición revisada)
$ws.Range("A5").Value = "DYmUGGwZ8_oC"
$ws.Range("B5").Value = "El Señor de los Anillos no 01/03 La Comunidad del Anillo (edición revisada)"
$ws.Range("C5").Value = "Primera entrega de la trilogía. «Este libro es como un relámpago en un cielo claro. Decir que la novela heroica, espléndida, elocuente y desinhibida, ha retornado de pronto en una época de un antirromanticismo casi patológico, sería inadecuado. Para quienes vivimos en esa extraña época, el retorno —y el alivio que nos trae— es sin duda lo más importante. Pero para la historia misma de la novela —una historia que se remonta a la Odisea y a antes de la Odisea— no es un retorno, sino un paso adelante o una revolución: la conquista de un territorio nuevo.» —C.S. Lewis, Time & Tide, 1954 «La obra de Tolkien, difundida en millones de ejemplares, traducida a docenas de lenguas, inspiradora de slogans pintados en las paredes de Nueva York y de Buenos Aires... una coherente mitología de una autenticidad universal creada en pleno siglo veinte.» —George Steiner, Le Monde, 1973"
$ws.Range("D5").Value = "2010-07-15"
$ws.Range("E5").Value = "J. R. R. Tolkien"

# Row 6: brand new content (Hobbit's Travels journal by Sam Gamgee)
$ws.Range("A6").Value = "FlGcUAnpMmIC"
$ws.Range("B6").Value = "Hobbit's Travels"
$ws.Range("C6").Value = "Printed on deluxe recycled parchment paper, this journal celebrating J. R. R. Tolkien's classic tales makes a lovely gift, and is just as nice to keep! With magical two-color illustrations throughout (drawings made by Frodo Baggins's devoted companion, Sam Gamgee, on their travels throughout Middle-earth), it provides ample space for recording personal thoughts, reflections on Tolkien's masterpiece, or fantasies of your own creation."
$ws.Range("D6").Value = "2002-01-04"
$ws.Range("E6").Value = "Sam Gamgee"
